# Adds three new "Use Case - Critical Infrastructure" slides to the end of
# the deck (slides 11, 12, 13), each using the "Title and Content" layout
# (layout index 2 -> slideLayout2.xml, same layout used by the rest of the
# deck's content slides).

$p = $ppt.ActivePresentation

$ppLayoutText = 2
$msoAutoSizeTextToFitShape = 2

$dash = [char]0x2013   # en dash "-"

# ---------------------------------------------------------------------
# Slide 11 - Use Case - Critical Infrastructure (overview)
# ---------------------------------------------------------------------
$s11 = $p.Slides.Add(11, $ppLayoutText)

$s11.Shapes.Item(1).TextFrame.TextRange.Text = "Use Case $dash Critical Infrastructure"

$body11 = $s11.Shapes.Item(2).TextFrame
$body11.AutoSize = $msoAutoSizeTextToFitShape
$tr11 = $body11.TextRange
$tr11.Text = "Target: A group of organizations that collaboratively manage critical infrastructure and utilize Industrial Control Systems.`rPower, water and other critical infrastructure are threatened by cyber and physical terrorism. `rIndustrial Control Systems are increasingly computer controlled and connected (directly or indirectly) to the internet and may embed compromised control hardware/software from questionable sources."

# ---------------------------------------------------------------------
# Slide 12 - Use Case - Critical Infrastructure (Scenario 1)
# ---------------------------------------------------------------------
$s12 = $p.Slides.Add(12, $ppLayoutText)

$s12.Shapes.Item(1).TextFrame.TextRange.Text = "Use Case $dash Critical Infrastructure"

$body12 = $s12.Shapes.Item(2).TextFrame
$body12.AutoSize = $msoAutoSizeTextToFitShape
$tr12 = $body12.TextRange
$tr12.Text = "Scenario 1 $dash North-East U.S. Power Grid attack by terrorists`rAn undetected and formerly unknown virus is planted in control systems emergency response software.`rA physical attack on a substation initiates a cascading failure`rCompromised control systems do not take substations and generating capacity off-line and introduce failure protocols, causing substantial failure of the physical infrastructure`rMuch of the power grid off-line for months"

for ($i = 2; $i -le 5; $i++) {
  $tr12.Paragraphs($i, 1).IndentLevel = 2
}

# ---------------------------------------------------------------------
# Slide 13 - Use Case - Critical Infrastructure (Scenario 1 Mitigations)
# ---------------------------------------------------------------------
$s13 = $p.Slides.Add(13, $ppLayoutText)

$s13.Shapes.Item(1).TextFrame.TextRange.Text = "Use Case $dash Critical Infrastructure"

$body13 = $s13.Shapes.Item(2).TextFrame
$body13.AutoSize = $msoAutoSizeTextToFitShape
$tr13 = $body13.TextRange
$tr13.Text = "Scenario 1 $dash Mitigations`rInformation from Control Systems, Field Monitors, Physical Security and User reports  are integrated via a threat management console`rInitial attack is recognized`rInitial Cascade is recognized`rFailure of control software generates incorrect action, which is recognized as an additional threat`rIntegrated threat management creates alert, manual control over infrastructure prevents widespread failure`r"

for ($i = 2; $i -le 7; $i++) {
  $tr13.Paragraphs($i, 1).IndentLevel = 2
}

# Split paragraph 1 into two runs: "Scenario 1 - " / "Mitigations"
$para1 = $tr13.Paragraphs(1, 1)
$firstLen = ("Scenario 1 $dash ").Length
$firstPart = $para1.Characters(1, $firstLen)
$firstPart.Text = $firstPart.Text

# Split paragraph 5 into two runs:
# "Failure ... recognized as an " / "additional threat"
$para5 = $tr13.Paragraphs(5, 1)
$prefix = "Failure of control software generates incorrect action, which is recognized as an "
$prefixPart = $para5.Characters(1, $prefix.Length)
$prefixPart.Text = $prefixPart.Text
